$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6251964.5
$ws.Range("J17").Value = 6668635.5
$ws.Range("L17").Value = 20005906.5
$ws.Range("N17").Value = -20006242.5
$ws.Range("H116").Value = 6027.091
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 6329.8
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 6329.8
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -13213.8
$ws.Range("H129").Value = 189675.94
$ws.Range("J129").Value = 201036.7
$ws.Range("L129").Value = 603110.1000000001
$ws.Range("N129").Value = -613110.1000000001
$ws.Range("H137").Value = 1011.64105
$ws.Range("I137").Value = 910.8
$ws.Range("K137").Value = 2732.4
$ws.Range("M137").Value = -182.3999999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 11503.5
$ws.Range("J27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("N27").Value = -12368
$ws.Range("H32").Value = 6755.9873
$ws.Range("I32").Value = 5511.9194
$ws.Range("J32").Value = 11293.177
$ws.Range("K32").Value = 5511.9194
$ws.Range("L32").Value = 11293.177
$ws.Range("M32").Value = -5224.9194
$ws.Range("N32").Value = -11867.177
$ws.Range("N33").ClearContents()
$ws.Range("H33").Value = 3666.6667
$ws.Range("I33").Value = 3666.6667
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 3666.6667
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3337.6667
$ws.Range("N34").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H39").Value = 3187
$ws.Range("I39").Value = 3187
$ws.Range("K39").Value = 3187
$ws.Range("M39").Value = -2667
$ws.Range("H45").Value = 3673
$ws.Range("I45").Value = 3416.75
$ws.Range("J45").Value = 3952.5454
$ws.Range("K45").Value = 3416.75
$ws.Range("L45").Value = 3952.5454
$ws.Range("M45").Value = -3039.75
$ws.Range("N45").Value = -4706.5454
$ws.Range("H61").Value = 2986.4482
$ws.Range("I61").Value = 2824.28
$ws.Range("K61").Value = 2824.28
$ws.Range("M61").Value = -2612.28
$ws.Range("H74").Value = 30304784
$ws.Range("I74").Value = 55556230
$ws.Range("J74").Value = 3047.4666
$ws.Range("K74").Value = 55556230
$ws.Range("L74").Value = 3047.4666
$ws.Range("M74").Value = -55555356
$ws.Range("N74").Value = -4795.4666
$ws.Range("H77").Value = 30304784
$ws.Range("I77").Value = 55556230
$ws.Range("J77").Value = 3047.4666
$ws.Range("K77").Value = 277781150
$ws.Range("L77").Value = 15237.333
$ws.Range("M77").Value = -277776782
$ws.Range("N77").Value = -23973.333
$ws.Range("H122").Value = 1346.7391
$ws.Range("I122").Value = 1371.591
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 4114.772999999999
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -1664.772999999999
$ws.Range("N122").Value = -7300
$ws.Range("N125").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H132").Value = 12290.667
$ws.Range("I132").Value = 1662.5641
$ws.Range("J132").Value = 58345.777
$ws.Range("K132").Value = 4987.692300000001
$ws.Range("L132").Value = 175037.331
$ws.Range("M132").Value = -2457.692300000001
$ws.Range("N132").Value = -180097.331
$ws.Range("H136").Value = 2986.4482
$ws.Range("I136").Value = 2824.28
$ws.Range("K136").Value = 8472.84
$ws.Range("M136").Value = -5922.84

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1677.75
$ws.Range("I99").Value = 1566.6666
$ws.Range("K99").Value = 1566.6666
$ws.Range("M99").Value = -68.66660000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3287.0227
$ws.Range("I31").Value = 1674.4445
$ws.Range("J31").Value = 5848.1763
$ws.Range("K31").Value = 1674.4445
$ws.Range("L31").Value = 5848.1763
$ws.Range("M31").Value = -1379.4445
$ws.Range("N31").Value = -6438.1763
$ws.Range("H34").Value = 3287.0227
$ws.Range("I34").Value = 1674.4445
$ws.Range("J34").Value = 5848.1763
$ws.Range("K34").Value = 1674.4445
$ws.Range("L34").Value = 5848.1763
$ws.Range("M34").Value = -1472.4445
$ws.Range("N34").Value = -6252.1763
$ws.Range("H122").Value = 1612.5
$ws.Range("I122").Value = 1483.3334
$ws.Range("K122").Value = 4450.0002
$ws.Range("M122").Value = -2000.0002
$ws.Range("H134").Value = 1118.3636
$ws.Range("I134").Value = 922.44446
$ws.Range("K134").Value = 2767.33338
$ws.Range("M134").Value = -232.33338

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 286.33334
$ws.Range("I8").Value = 286.33334
$ws.Range("K8").Value = 859.0000200000001
$ws.Range("M8").Value = -720.0000200000001
$ws.Range("H115").Value = 4680.4443
$ws.Range("J115").Value = 5261.75
$ws.Range("L115").Value = 15785.25
$ws.Range("N115").Value = -18135.25
$ws.Range("H131").Value = 690.14
$ws.Range("J131").Value = 717.73627
$ws.Range("L131").Value = 2153.20881
$ws.Range("N131").Value = -12233.20881
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 454.875
$ws.Range("I132").Value = 454.875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4093.875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1563.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3513.3333
$ws.Range("I80").Value = 3046.5386
$ws.Range("J80").Value = 3946.7856
$ws.Range("K80").Value = 3046.5386
$ws.Range("L80").Value = 3946.7856
$ws.Range("M80").Value = -2048.5386
$ws.Range("N80").Value = -5942.7856
$ws.Range("H83").Value = 3513.3333
$ws.Range("I83").Value = 3046.5386
$ws.Range("J83").Value = 3946.7856
$ws.Range("K83").Value = 15232.693
$ws.Range("L83").Value = 19733.928
$ws.Range("M83").Value = -10240.693
$ws.Range("N83").Value = -29717.928
$ws.Range("H97").Value = 1708.3334
$ws.Range("I97").Value = 1847.7273
$ws.Range("J97").Value = 1325
$ws.Range("K97").Value = 1847.7273
$ws.Range("L97").Value = 1325
$ws.Range("M97").Value = -1351.7273
$ws.Range("N97").Value = -2317
$ws.Range("H102").Value = 1885.931
$ws.Range("I102").Value = 1608.7693
$ws.Range("K102").Value = 1608.7693
$ws.Range("M102").Value = 13.23070000000007
$ws.Range("H122").Value = 4584.5
$ws.Range("I122").Value = 1835.6666
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 5506.9998
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -3056.9998
$ws.Range("N122").Value = -26900.0005
$ws.Range("H123").Value = 8656.625
$ws.Range("I123").Value = 4620
$ws.Range("J123").Value = 10002.167
$ws.Range("K123").Value = 4620
$ws.Range("L123").Value = 10002.167
$ws.Range("M123").Value = -2170
$ws.Range("N123").Value = -14902.167
$ws.Range("H132").Value = 25403.875
$ws.Range("I132").Value = 4935.294
$ws.Range("J132").Value = 75113.28999999999
$ws.Range("K132").Value = 14805.882
$ws.Range("L132").Value = 225339.87
$ws.Range("M132").Value = -12275.882
$ws.Range("N132").Value = -230399.87

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5444.077
$ws.Range("I61").Value = 2086
$ws.Range("K61").Value = 2086
$ws.Range("M61").Value = -1884
$ws.Range("H113").Value = 5444.077
$ws.Range("I113").Value = 2086
$ws.Range("K113").Value = 2086
$ws.Range("M113").Value = 84
$ws.Range("H122").Value = 1035849.8
$ws.Range("I122").Value = 1403703.4
$ws.Range("K122").Value = 4211110.199999999
$ws.Range("M122").Value = -4208660.199999999
$ws.Range("N123").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H136").Value = 1685.9546
$ws.Range("I136").Value = 1480.5238
$ws.Range("K136").Value = 4441.5714
$ws.Range("M136").Value = -1891.5714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1796.5714
$ws.Range("I122").Value = 1775.2
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 5325.6
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -2875.6
$ws.Range("N122").Value = -10450
$ws.Range("N123").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H125").Value = 48999
$ws.Range("J125").Value = 48999
$ws.Range("L125").Value = 48999
$ws.Range("N125").Value = -58839
$ws.Range("H132").Value = 1219.5319
$ws.Range("I132").Value = 843.6111
$ws.Range("J132").Value = 2449.818
$ws.Range("K132").Value = 2530.8333
$ws.Range("L132").Value = 7349.454000000001
$ws.Range("M132").Value = -0.833299999999781
$ws.Range("N132").Value = -12409.454
